$d = $word.ActiveDocument
$VT = [char]11   # vertical-tab char = how Word represents a manual <w:br/> in Range.Text

# ---------------------------------------------------------------------------
# CHANGE 1: paragraph 1.2 rewording
# ---------------------------------------------------------------------------

# 1a) drop the "/ специальности (выбрать нужное)" choice-prompt
$d.Content.Find.Execute(
    "по направлению подготовки / специальности (выбрать нужное) 09",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "по направлению подготовки 09", 2) | Out-Null

# 1b) drop "производственной научной и преддипломной практики" + the manual
#     line break that followed "места для прохождения "
$d.Content.Find.Execute(
    "места для прохождения производственной научной и преддипломной практики $VT",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "места для прохождения ", 2) | Out-Null

# 1c) insert the new bold+underlined "производственной практики" run, followed
#     by an underlined parenthetical clarification, right before "обучающихся"
$r = $d.Content
$r.Find.Execute("места для прохождения ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0) | Out-Null                      # wdCollapseEnd

$r.InsertAfter("производственной практики")
$r.Font.Size = 8
$r.Font.Bold = 1
$r.Font.Underline = 1

$r.Collapse(0) | Out-Null
$r.InsertAfter(" (практики по получению профессиональных умений и опыта профессиональной деятельности, научно-исследовательская работы и преддипломной практики)")
$r.Font.Size = 8
$r.Font.Bold = 0
$r.Font.Underline = 1

# ---------------------------------------------------------------------------
# CHANGE 2: «5.3.» contract-validity date
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "твует с «06» 04 2021 г. по «____»____________ 20__ г.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "твует с «05» 04 2021 г. по «____»____________ 2021 г.", 2) | Out-Null

# NOTE: the source revision also shifts a <w:lastRenderedPageBreak/> marker
# from clause 3.4 to clause 3.3. That element is a purely internal,
# auto-computed Word pagination cache hint (it has no visible/content
# effect and is not exposed as a settable property anywhere on the Word
# object model - real Word VBA cannot place it either), so it is
# intentionally left untouched here.

Write-Host "edit.ps1 complete"
